$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the header cells in row 1: "_old" suffix group -> "_FV2404", "_new" suffix group -> "_FV2410"
$headersFV2404 = @("Segmentname_FV2404","Segmentgruppe_FV2404","Segment_FV2404","Datenelement_FV2404","Segment ID_FV2404","Code_FV2404","Qualifier_FV2404","Beschreibung_FV2404","Bedingungsausdruck_FV2404","Bedingung_FV2404")
for ($i = 0; $i -lt $headersFV2404.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headersFV2404[$i]
}

$headersFV2410 = @("Segmentname_FV2410","Segmentgruppe_FV2410","Segment_FV2410","Datenelement_FV2410","Segment ID_FV2410","Code_FV2410","Qualifier_FV2410","Beschreibung_FV2410","Bedingungsausdruck_FV2410","Bedingung_FV2410")
for ($i = 0; $i -lt $headersFV2410.Length; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $headersFV2410[$i]
}

# Turn the data range into an Excel Table ("Table1") with the header row.
$tableRange = $ws.Range("A1:U86")
$listObj = $ws.ListObjects.Add(1, $tableRange, 0, 1)
$listObj.Name = "Table1"
$listObj.TableStyle = ""

# Freeze the header row (first row stays visible while scrolling).
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
